$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> (predidx, pred_name)
$updates = @(
    @{ Row = 3;   D = "[1, 0, 0, 1, 0, 0, 0]"; E = "['Normal', 'ParamViolation']" },
    @{ Row = 15;  D = "[0, 0, 0, 1, 0, 0, 0]"; E = "['ParamViolation']" },
    @{ Row = 25;  D = "[1, 0, 0, 0, 0, 0, 1]"; E = "['Normal', 'SoftwareFault']" },
    @{ Row = 38;  D = "[0, 0, 1, 0, 0, 0, 1]"; E = "['HardwareFault', 'SoftwareFault']" },
    @{ Row = 39;  D = "[1, 0, 0, 0, 0, 0, 1]"; E = "['Normal', 'SoftwareFault']" },
    @{ Row = 54;  D = "[0, 0, 0, 0, 0, 0, 1]"; E = "['SoftwareFault']" },
    @{ Row = 56;  D = "[0, 0, 0, 0, 0, 0, 0]"; E = "[]" },
    @{ Row = 61;  D = "[1, 0, 0, 0, 0, 0, 1]"; E = "['Normal', 'SoftwareFault']" },
    @{ Row = 68;  D = "[1, 0, 0, 0, 0, 0, 0]"; E = "['Normal']" },
    @{ Row = 69;  D = "[1, 1, 0, 0, 0, 0, 0]"; E = "['Normal', 'SurroundingEnvironment']" },
    @{ Row = 70;  D = "[1, 1, 0, 0, 0, 0, 0]"; E = "['Normal', 'SurroundingEnvironment']" },
    @{ Row = 73;  D = "[1, 0, 0, 0, 0, 0, 0]"; E = "['Normal']" },
    @{ Row = 82;  D = "[1, 1, 1, 0, 0, 0, 0]"; E = "['Normal', 'SurroundingEnvironment', 'HardwareFault']" },
    @{ Row = 88;  D = "[1, 0, 0, 0, 0, 0, 0]"; E = "['Normal']" },
    @{ Row = 109; D = "[1, 0, 0, 0, 0, 0, 1]"; E = "['Normal', 'SoftwareFault']" },
    @{ Row = 113; D = "[1, 0, 0, 0, 0, 0, 1]"; E = "['Normal', 'SoftwareFault']" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 4).Value = $u.D
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
